$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 573, pushing the existing rows 573-661 down to
# 576-664, to make room for a new weekly price entry.
$ws.Rows("573:575").Insert()

# Common (constant) field values shared by every data row in this sheet.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100106
$producto    = "Oleaginosos"
$categoriaId = 100106002
$categoria   = "Palta"
$origen      = "Provincia de Limarí"
$kgUnidad    = 1

# New weekly entry: Comercializadora del Agro de Limarí - Palta - Hass,
# fecha 2023-03-30 (serial 45015), caja de 17 kilos.
$dateVal = 45015
$unidad  = "$/kilo (en caja de 17 kilos)"

$newRows = @(
    @{ Row = 573; Calidad = "Especial"; Volumen = 500; Min = 4000; Max = 4100; Prom = 4050 },
    @{ Row = 574; Calidad = "Primera";  Volumen = 300; Min = 3800; Max = 3900; Prom = 3850 },
    @{ Row = 575; Calidad = "Segunda";  Volumen = 240; Min = 3300; Max = 3400; Prom = 3350 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $dateVal
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $tipo
    $ws.Range("G$row").Value = $productoId
    $ws.Range("H$row").Value = $producto
    $ws.Range("I$row").Value = $categoriaId
    $ws.Range("J$row").Value = $categoria
    $ws.Range("K$row").Value = "Hass"
    $ws.Range("L$row").Value = $r.Calidad
    $ws.Range("M$row").Value = $r.Volumen
    $ws.Range("N$row").Value = $r.Min
    $ws.Range("O$row").Value = $r.Max
    $ws.Range("P$row").Value = $r.Prom
    $ws.Range("Q$row").Value = $unidad
    $ws.Range("R$row").Value = $origen
    $ws.Range("S$row").Value = $r.Prom
    $ws.Range("T$row").Value = $kgUnidad
}
